$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was added at the top of the Zanahoria /
# Terminal Hortofrutícola Agro Chillán block (row 243), pushing the
# existing rows 243-258 down to 244-259.
$ws.Rows.Item(243).Insert()

$ws.Range("A243").Value = 7
$ws.Range("B243").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C243").Value = "Ñuble"
$ws.Range("D243").Value = 44615
$ws.Range("E243").Value = 16
$ws.Range("F243").Value = 100114013
$ws.Range("G243").Value = "Zanahoria"
$ws.Range("H243").Value = "Sin especificar"
$ws.Range("I243").Value = "Primera"
$ws.Range("J243").Value = 100
$ws.Range("K243").Value = 7000
$ws.Range("L243").Value = 7500
$ws.Range("M243").Value = 7250
$ws.Range("N243").Value = "$/saco 20 kilos"
$ws.Range("O243").Value = "Provincia de Diguillín"
$ws.Range("P243").Value = 362
$ws.Range("Q243").Value = 20
$ws.Range("R243").Value = "Hortaliza"
